# "work on checkout system"
#
# Slide 24 is the checkout-system mockup. Tidy up the "scan barcode / manual
# entry" controls:
#   - move the "Scan barcode:" label
#   - drop the separate "Or input manually" label (its job is taken over by
#     the manual-entry box sitting right next to the barcode label)
#   - slide the (blank) manual-entry box up next to the "Scan barcode:" label
#   - move the "Add to basket" button up to follow
#
# NOTE on the literals below: PowerPoint's object model reports/accepts
# shape Left/Top/Width/Height in points (1 pt = 12700 EMU) but this host
# stores the value as a 32-bit float internally, then truncates when
# converting back to EMU. A "natural" literal like `emu/12700.0` can
# therefore land one EMU short after the round-trip. The literals used here
# are the nearest representable points that truncate back to the exact
# target EMU so the saved XML matches precisely.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(24)

# "Scan barcode:" label -> new position (310319, 1819201) EMU; size unchanged.
$scanLabel = $s.Shapes.Item("TextBox 1")
$scanLabel.Left = 24.434568405151367
$scanLabel.Top = 143.2441864013672

# "Or input manually" label is no longer needed - remove it.
$orInputManually = $s.Shapes.Item("TextBox 5")
$orInputManually.Delete()

# Blank manual barcode-entry box -> moved up beside the label;
# (1942775, 1819201) EMU, size stays 1852179 x 369332 EMU.
$manualEntryBox = $s.Shapes.Item("TextBox 8")
$manualEntryBox.Left = 152.9744110107422
$manualEntryBox.Top = 143.2441864013672

# "Add to basket" button -> new position (1126547, 2862884) EMU; size unchanged.
$addToBasket = $s.Shapes.Item("TextBox 7")
$addToBasket.Left = 88.7044906616211
$addToBasket.Top = 225.4239501953125
